$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("mvmt -3 pips").Name = "-3 pips"
$wb.Worksheets.Item("mvmt -4 pips").Name = "-4 pips"
$wb.Worksheets.Item("mvmt -5 pips").Name = "-5 pips"
$wb.Worksheets.Item("mvmt -6 pips").Name = "-6 pips"
